$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 2 (hSSQ1lIH - Athletico-PR vs Atletico-MG) and row 3 (2BPmIXwm - Amazonas vs Goias)
$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.25
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("AA2").Value = 21
$ws.Range("AH2").Value = 8.5
$ws.Range("AN2").Value = 4.33
$ws.Range("G3").Value = 3.6
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 2.05
$ws.Range("J3").Value = 4.33
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 2.88
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("W3").Value = 9.5
$ws.Range("X3").Value = 17
$ws.Range("AC3").Value = 8.5
$ws.Range("AK3").Value = 19
$ws.Range("AL3").Value = 19
$ws.Range("AM3").Value = 34
$ws.Range("AQ3").Value = 67
$ws.Range("AT3").Value = 2.5
$ws.Range("AV3").Value = 67
$ws.Range("AX3").Value = 12
$ws.Range("BA3").Value = 67
$ws.Range("BB3").Value = 201

# Remove row 6 (pCxQFr9d - Sp. Luqueno vs Ameliano); subsequent rows shift up
$ws.Rows(6).Delete()
